$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 81) to the daily log sheet:
# 2025/10/09, 木 (Thursday), time code 1, ranking 125

# Force column A to be treated as text so the date-like string
# "2025/10/09" is stored as a literal string (matching the existing
# rows) instead of being auto-converted into an Excel date serial
# number, then restore the cell's style to the default "Normal"
# style so no extra formatting is left behind on the cell.
$ws.Range("A81").NumberFormat = "@"
$ws.Range("A81").Value = "2025/10/09"
$ws.Range("A81").Style = "Normal"

$ws.Range("B81").Value = "木"
$ws.Range("C81").Value = 1
$ws.Range("D81").Value = 125
